$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "44.060.70"
Set-TextValue $ws.Range("E2") "  +0.49%  "
Set-TextValue $ws.Range("D3") "2.328.15"
Set-TextValue $ws.Range("E3") "  +4.11%  "
Set-TextValue $ws.Range("E4") "  +0.07%  "
Set-TextValue $ws.Range("D5") "98.49"
Set-TextValue $ws.Range("E5") "  +4.44%  "
Set-TextValue $ws.Range("D6") "271.53"
Set-TextValue $ws.Range("E6") "  +0.02%  "
Set-TextValue $ws.Range("E7") "  +0.28%  "
Set-TextValue $ws.Range("E8") "  +0.00%  "
Set-TextValue $ws.Range("D9") "0.629"
Set-TextValue $ws.Range("E9") "  -1.13%  "
Set-TextValue $ws.Range("E10") "  -1.17%  "
Set-TextValue $ws.Range("D11") "0.0958"
Set-TextValue $ws.Range("E11") "  +0.27%  "
Set-TextValue $ws.Range("D12") "8.04"
Set-TextValue $ws.Range("E12") "  -3.98%  "
Set-TextValue $ws.Range("E13") "  +0.41%  "
Set-TextValue $ws.Range("D14") "2.665.55"
Set-TextValue $ws.Range("E14") "  +3.63%  "
Set-TextValue $ws.Range("D15") "15.58"
Set-TextValue $ws.Range("E15") "  +1.40%  "
Set-TextValue $ws.Range("D16") "0.882"
Set-TextValue $ws.Range("E16") "  +7.02%  "
Set-TextValue $ws.Range("D17") "2.325.19"
Set-TextValue $ws.Range("D18") "44.012.46"
Set-TextValue $ws.Range("E19") "  +4.57%  "
Set-TextValue $ws.Range("D20") "6.40"
Set-TextValue $ws.Range("E20") "  +3.69%  "
Set-TextValue $ws.Range("D21") "73.73"
Set-TextValue $ws.Range("E21") "  +4.04%  "
Set-TextValue $ws.Range("D22") "2.32"
Set-TextValue $ws.Range("E22") "  -0.85%  "
Set-TextValue $ws.Range("D23") "240.44"
Set-TextValue $ws.Range("E23") "  +2.46%  "
Set-TextValue $ws.Range("D24") "9.29"
Set-TextValue $ws.Range("E24") "  +1.32%  "
Set-TextValue $ws.Range("E25") "  -0.05%  "
Set-TextValue $ws.Range("E26") "  +1.26%  "
Set-TextValue $ws.Range("E27") "  -0.26%  "
Set-TextValue $ws.Range("E28") "  -1.89%  "
Set-TextValue $ws.Range("E29") "  +1.94%  "
Set-TextValue $ws.Range("D30") "38.50"
Set-TextValue $ws.Range("E30") "  -4.53%  "
Set-TextValue $ws.Range("D31") "22.49"
Set-TextValue $ws.Range("E31") "  +7.06%  "
Set-TextValue $ws.Range("D32") "176.05"
Set-TextValue $ws.Range("E32") "  +1.96%  "
Set-TextValue $ws.Range("D33") "0.0917"
Set-TextValue $ws.Range("E33") "  -0.01%  "
Set-TextValue $ws.Range("E34") "  +0.42%  "
Set-TextValue $ws.Range("E35") "  +2.06%  "
Set-TextValue $ws.Range("D36") "0.0365"
Set-TextValue $ws.Range("E36") "  +3.53%  "
Set-TextValue $ws.Range("E37") "  -2.35%  "
Set-TextValue $ws.Range("E38") "  +3.98%  "
Set-TextValue $ws.Range("D39") "3.39"
Set-TextValue $ws.Range("E39") "  -4.69%  "
Set-TextValue $ws.Range("D40") "2.44"
Set-TextValue $ws.Range("E40") "  +13.41%  "
Set-TextValue $ws.Range("E41") "  +9.06%  "
Set-TextValue $ws.Range("E42") "  +25.36%  "
Set-TextValue $ws.Range("D43") "12.47"
Set-TextValue $ws.Range("E43") "  -2.91%  "
Set-TextValue $ws.Range("D44") "62.80"
Set-TextValue $ws.Range("E44") "  -0.84%  "
Set-TextValue $ws.Range("E45") "  +8.83%  "
Set-TextValue $ws.Range("D46") "5.34"
Set-TextValue $ws.Range("E46") "  -1.26%  "
Set-TextValue $ws.Range("E47") "  +4.18%  "
Set-TextValue $ws.Range("D48") "100.58"
Set-TextValue $ws.Range("E48") "  -1.25%  "
Set-TextValue $ws.Range("E49") "  +0.44%  "
Set-TextValue $ws.Range("E50") "  +16.21%  "
Set-TextValue $ws.Range("D51") "2.552.50"
Set-TextValue $ws.Range("E51") "  +3.96%  "
